# Commit: "changed stock to holding"
# Rename the "Stock" and "Balance" column headers to "Total Holding" and
# "Total Balance" respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Total Holding"
$ws.Range("D1").Value = "Total Balance"

# Give the bordered box around the "Init" row a solid white fill so it
# reads as an opaque panel instead of a transparent one (matches the
# existing header row's fill treatment).
$ws.Range("A2:F2").Interior.ColorIndex = 2
$ws.Range("A3").Interior.ColorIndex = 2
